$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell B7 from "str" to "list.str"
$ws.Range("B7").Value = "list.str"

# Update the selection to B8 (as reflected in sheetView)
$ws.Range("B8").Select()
